$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename PAN header to Pan
$ws.Range("H1").Value = "Pan"

# Add new headers for DP and Client Id columns
$ws.Range("N1").Value = "DP"
$ws.Range("O1").Value = "Client Id"

# Fill in DP (N) and Client Id (O) values for each data row
$ws.Range("N2").Value = 123456
$ws.Range("O2").Value = 1

$ws.Range("N3").Value = 234567
$ws.Range("O3").Value = 2

$ws.Range("N4").Value = 345678
$ws.Range("O4").Value = 3

$ws.Range("N5").Value = 456789
$ws.Range("O5").Value = 4

$ws.Range("N6").Value = 567900
$ws.Range("O6").Value = 5

$ws.Range("N7").Value = 679011
$ws.Range("O7").Value = 6

# Match the new styling used for header cells M1/N1/O1 (style index 4, same as M1)
$ws.Range("N1:O1").Style = $ws.Range("M1").Style

# Update the sheet view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("O8").Select() | Out-Null
